$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3595.8333
$ws.Range("I62").Value = 2894.8572
$ws.Range("K62").Value = 2894.8572
$ws.Range("M62").Value = -2270.8572
$ws.Range("H65").Value = 3595.8333
$ws.Range("I65").Value = 2894.8572
$ws.Range("K65").Value = 14474.286
$ws.Range("M65").Value = -11354.286
$ws.Range("H107").Value = 1567.45
$ws.Range("I107").Value = 1597.1428
$ws.Range("K107").Value = 1597.1428
$ws.Range("M107").Value = 322.8571999999999
$ws.Range("H111").Value = 966.6667
$ws.Range("I111").Value = 966.6667
$ws.Range("K111").Value = 2900.0001
$ws.Range("M111").Value = 166.9998999999998
$ws.Range("H112").Value = 1236.5646
$ws.Range("I112").Value = 350
$ws.Range("J112").Value = 1297.7069
$ws.Range("K112").Value = 1050
$ws.Range("L112").Value = 3893.120699999999
$ws.Range("M112").Value = 58
$ws.Range("N112").Value = -6109.120699999999
$ws.Range("H132").Value = 32260076
$ws.Range("I132").Value = 34484324
$ws.Range("J132").Value = 8500
$ws.Range("K132").Value = 103452972
$ws.Range("L132").Value = 25500
$ws.Range("M132").Value = -103450442
$ws.Range("N132").Value = -30560
$ws.Range("H137").Value = 4765147.5
$ws.Range("I137").Value = 9525395
$ws.Range("J137").Value = 4900
$ws.Range("K137").Value = 28576185
$ws.Range("L137").Value = 14700
$ws.Range("M137").Value = -28573635
$ws.Range("N137").Value = -19800
$ws.Range("H138").Value = 2660.49
$ws.Range("I138").Value = 793.93335
$ws.Range("J138").Value = 2989.8823
$ws.Range("K138").Value = 2381.80005
$ws.Range("L138").Value = 8969.6469
$ws.Range("M138").Value = 2758.19995
$ws.Range("N138").Value = -19249.6469

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7066.9165
$ws.Range("I32").Value = 5803.516
$ws.Range("J32").Value = 14900
$ws.Range("K32").Value = 5803.516
$ws.Range("L32").Value = 14900
$ws.Range("M32").Value = -5516.516
$ws.Range("N32").Value = -15474
$ws.Range("H61").Value = 1585.6923
$ws.Range("I61").Value = 1320
$ws.Range("J61").Value = 2471.3333
$ws.Range("K61").Value = 1320
$ws.Range("L61").Value = 2471.3333
$ws.Range("M61").Value = -1108
$ws.Range("N61").Value = -2895.3333
$ws.Range("H74").Value = 6418.8335
$ws.Range("I74").Value = 8242.362999999999
$ws.Range("K74").Value = 8242.362999999999
$ws.Range("M74").Value = -7368.362999999999
$ws.Range("H77").Value = 6418.8335
$ws.Range("I77").Value = 8242.362999999999
$ws.Range("K77").Value = 41211.815
$ws.Range("M77").Value = -36843.815
$ws.Range("H124").Value = 29095.334
$ws.Range("J124").Value = 29095.334
$ws.Range("L124").Value = 29095.334
$ws.Range("N124").Value = -38915.334
$ws.Range("H132").Value = 4687.5
$ws.Range("I132").Value = 1712
$ws.Range("J132").Value = 5282.6
$ws.Range("K132").Value = 5136
$ws.Range("L132").Value = 15847.8
$ws.Range("M132").Value = -2606
$ws.Range("N132").Value = -20907.8
$ws.Range("H136").Value = 1585.6923
$ws.Range("I136").Value = 1320
$ws.Range("J136").Value = 2471.3333
$ws.Range("K136").Value = 3960
$ws.Range("L136").Value = 7413.999899999999
$ws.Range("M136").Value = -1410
$ws.Range("N136").Value = -12513.9999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 16363.909
$ws.Range("I20").Value = 2750.5
$ws.Range("J20").Value = 32700
$ws.Range("K20").Value = 2750.5
$ws.Range("L20").Value = 32700
$ws.Range("M20").Value = -2503.5
$ws.Range("N20").Value = -33194
$ws.Range("H62").Value = 43500
$ws.Range("J62").Value = 43500
$ws.Range("L62").Value = 43500
$ws.Range("N62").Value = -44872
$ws.Range("H65").Value = 43500
$ws.Range("J65").Value = 43500
$ws.Range("L65").Value = 130500
$ws.Range("N65").Value = -137364
$ws.Range("H134").Value = 3086.1304
$ws.Range("I134").Value = 2428.2942
$ws.Range("K134").Value = 7284.882599999999
$ws.Range("M134").Value = -4749.882599999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2909.3125
$ws.Range("I31").Value = 1131.1111
$ws.Range("J31").Value = 5195.5713
$ws.Range("K31").Value = 1131.1111
$ws.Range("L31").Value = 5195.5713
$ws.Range("M31").Value = -836.1111000000001
$ws.Range("N31").Value = -5785.5713
$ws.Range("H34").Value = 2909.3125
$ws.Range("I34").Value = 1131.1111
$ws.Range("J34").Value = 5195.5713
$ws.Range("K34").Value = 1131.1111
$ws.Range("L34").Value = 5195.5713
$ws.Range("M34").Value = -929.1111000000001
$ws.Range("N34").Value = -5599.5713
$ws.Range("H58").Value = 3080.4688
$ws.Range("I58").Value = 1767.6666
$ws.Range("J58").Value = 8230.691999999999
$ws.Range("K58").Value = 1767.6666
$ws.Range("L58").Value = 8230.691999999999
$ws.Range("M58").Value = -1564.6666
$ws.Range("N58").Value = -8636.691999999999
$ws.Range("H68").Value = 92250.42999999999
$ws.Range("J68").Value = 92250.42999999999
$ws.Range("L68").Value = 92250.42999999999
$ws.Range("N68").Value = -93748.42999999999
$ws.Range("H71").Value = 92250.42999999999
$ws.Range("J71").Value = 92250.42999999999
$ws.Range("L71").Value = 276751.29
$ws.Range("N71").Value = -284239.29
$ws.Range("H134").Value = 2482.4119
$ws.Range("I134").Value = 1548.9166
$ws.Range("J134").Value = 4722.8
$ws.Range("K134").Value = 4646.7498
$ws.Range("L134").Value = 14168.4
$ws.Range("M134").Value = -2111.7498
$ws.Range("N134").Value = -19238.4
$ws.Range("H136").Value = 3080.4688
$ws.Range("I136").Value = 1767.6666
$ws.Range("J136").Value = 8230.691999999999
$ws.Range("K136").Value = 5302.9998
$ws.Range("L136").Value = 24692.076
$ws.Range("M136").Value = -2752.9998
$ws.Range("N136").Value = -29792.076
$ws.Range("H137").Value = 43927.145
$ws.Range("J137").Value = 43927.145
$ws.Range("L137").Value = 43927.145
$ws.Range("N137").Value = -54127.145
$ws.Range("H138").Value = 42722.453
$ws.Range("J138").Value = 42722.453
$ws.Range("L138").Value = 42722.453
$ws.Range("N138").Value = -53002.453
$ws.Range("H140").Value = 83793.84
$ws.Range("J140").Value = 83793.84
$ws.Range("L140").Value = 83793.84
$ws.Range("N140").Value = -94153.84

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1972.409
$ws.Range("I129").Value = 2014.8182
$ws.Range("J129").Value = 1930
$ws.Range("K129").Value = 6044.4546
$ws.Range("L129").Value = 5790
$ws.Range("M129").Value = -1044.4546
$ws.Range("N129").Value = -15790
$ws.Range("H131").Value = 692.5
$ws.Range("I131").Value = 258.1
$ws.Range("J131").Value = 806.8158
$ws.Range("K131").Value = 774.3000000000001
$ws.Range("L131").Value = 2420.4474
$ws.Range("M131").Value = 4265.7
$ws.Range("N131").Value = -12500.4474

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3364.8462
$ws.Range("I102").Value = 2245.8572
$ws.Range("J102").Value = 4670.3335
$ws.Range("K102").Value = 2245.8572
$ws.Range("L102").Value = 4670.3335
$ws.Range("M102").Value = -623.8571999999999
$ws.Range("N102").Value = -7914.3335
$ws.Range("H107").Value = 13889730
$ws.Range("I107").Value = 396
$ws.Range("J107").Value = 18519508
$ws.Range("K107").Value = 396
$ws.Range("L107").Value = 18519508
$ws.Range("M107").Value = 1524
$ws.Range("N107").Value = -18523348
$ws.Range("H132").Value = 3681.4736
$ws.Range("I132").Value = 1492
$ws.Range("J132").Value = 5273.8184
$ws.Range("K132").Value = 4476
$ws.Range("L132").Value = 15821.4552
$ws.Range("M132").Value = -1946
$ws.Range("N132").Value = -20881.4552

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 663.4607
$ws.Range("I68").Value = 663.4607
$ws.Range("K68").Value = 663.4607
$ws.Range("M68").Value = 85.53930000000003
$ws.Range("H71").Value = 663.4607
$ws.Range("I71").Value = 663.4607
$ws.Range("K71").Value = 3317.3035
$ws.Range("M71").Value = 426.6965
$ws.Range("H93").Value = 9262315
$ws.Range("I93").Value = 18521368
$ws.Range("J93").Value = 3262.6667
$ws.Range("K93").Value = 18521368
$ws.Range("L93").Value = 3262.6667
$ws.Range("M93").Value = -18520120
$ws.Range("N93").Value = -5758.6667

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31622
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -98112
$ws.Range("H107").Value = 622.5833
$ws.Range("I107").Value = 550.9375
$ws.Range("J107").Value = 765.875
$ws.Range("K107").Value = 1652.8125
$ws.Range("L107").Value = 2297.625
$ws.Range("M107").Value = 267.1875
$ws.Range("N107").Value = -6137.625
$ws.Range("H138").Value = 37649.668
$ws.Range("J138").Value = 37649.668
$ws.Range("L138").Value = 37649.668
$ws.Range("N138").Value = -47929.668

Write-Output "Applied 227 cell updates."